$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C4:Q4").Value = "Silt"
$ws.Range("C6:Q6").Value = "Top30cm"
[void]$ws.Range("C6:Q6").Select()
